# Refresh the stock-screener figures on Sheet1 (current price, change %,
# Stochastic Fast %K/%D) with the updated source data. Values below are
# written as plain decimal literals (not scientific notation) because the
# interpreter here does not accept an "E" exponent in numeric literals;
# Excel itself re-serializes them back to the compact XML form on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24150
$ws.Range("D2").Value = 0.0020999999999999999
$ws.Range("I2").Value = 4.1399999999999997
$ws.Range("J2").Value = 84
$ws.Range("K2").Value = 84

$ws.Range("C3").Value = 102300
$ws.Range("D3").Value = -0.073400000000000007
$ws.Range("I3").Value = 6.35
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 66

$ws.Range("C4").Value = 440500
$ws.Range("D4").Value = -0.0079000000000000008
$ws.Range("I4").Value = 4.3099999999999996

$ws.Range("C5").Value = 32150
$ws.Range("D5").Value = -0.0045999999999999999
$ws.Range("I5").Value = 6.22
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 50

$ws.Range("C6").Value = 31350
$ws.Range("D6").Value = 0.0016000000000000001
$ws.Range("D6").NumberFormat = "0.00%"
$ws.Range("I6").Value = 3.83
$ws.Range("J6").Value = 81
$ws.Range("K6").Value = 81

$ws.Range("C7").Value = 25550
$ws.Range("I7").Value = 4.7

$ws.Range("C8").Value = 10660
$ws.Range("D8").Value = 0.0057000000000000002
$ws.Range("I8").Value = 4.83
$ws.Range("J8").Value = 86
$ws.Range("K8").Value = 86

$ws.Range("C9").Value = 87000
$ws.Range("D9").Value = -0.0033999999999999998
$ws.Range("I9").Value = 3.45
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 80

$ws.Range("C10").Value = 213000
$ws.Range("D10").Value = -0.0448
$ws.Range("I10").Value = 5.63
$ws.Range("J10").Value = 51
$ws.Range("K10").Value = 51

$ws.Range("C11").Value = 128000
$ws.Range("D11").Value = -0.013899999999999999
$ws.Range("I11").Value = 5.31
$ws.Range("J11").Value = 82
$ws.Range("K11").Value = 82

$ws.Range("C12").Value = 20200
$ws.Range("D12").Value = 0.01
$ws.Range("D12").NumberFormat = "0.00%"
$ws.Range("I12").Value = 4.7

$ws.Range("C13").Value = 70800
$ws.Range("D13").Value = -0.0028
$ws.Range("I13").Value = 4.9400000000000004
$ws.Range("J13").Value = 81
$ws.Range("K13").Value = 81

$ws.Range("C14").Value = 56100
$ws.Range("D14").Value = 0.0018
$ws.Range("I14").Value = 6.31
$ws.Range("J14").Value = 73
$ws.Range("K14").Value = 73

$ws.Range("C15").Value = 85300
$ws.Range("D15").Value = 0.054399999999999997
$ws.Range("I15").Value = 6.45
$ws.Range("J15").Value = 91
$ws.Range("K15").Value = 91

$ws.Range("C16").Value = 19420
$ws.Range("D16").Value = -0.0041000000000000003
$ws.Range("I16").Value = 5.48
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 80

$ws.Range("C17").Value = 49800
$ws.Range("D17").Value = 0.001
$ws.Range("I17").Value = 5.62

$ws.Range("C18").Value = 19880
$ws.Range("D18").Value = -0.0085000000000000006
$ws.Range("I18").Value = 6.19
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 34

$ws.Range("C19").Value = 55200
$ws.Range("D19").Value = 0.010999999999999999
$ws.Range("I19").Value = 3.62
$ws.Range("J19").Value = 89
$ws.Range("K19").Value = 89

$ws.Range("C20").Value = 14630
$ws.Range("D20").Value = 0
$ws.Range("D20").NumberFormat = "0%"
$ws.Range("I20").Value = 4.4400000000000004
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 77

$ws.Range("C21").Value = 130700
$ws.Range("D21").Value = -0.0061000000000000004

$ws.Range("C22").Value = 42950
$ws.Range("D22").Value = -0.034799999999999998
$ws.Range("I22").Value = 3.39
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = 48

$ws.Range("C23").Value = 68000
$ws.Range("D23").Value = -0.013100000000000001
$ws.Range("I23").Value = 3.18
$ws.Range("J23").Value = 88
$ws.Range("K23").Value = 88

$ws.Range("C24").Value = 49300
$ws.Range("I24").Value = 5.48

$ws.Range("C25").Value = 85400
$ws.Range("D25").Value = -0.0070000000000000001
$ws.Range("I25").Value = 4.22
$ws.Range("J25").Value = 83
$ws.Range("K25").Value = 83

$ws.Range("C26").Value = 110900
$ws.Range("D26").Value = -0.016799999999999999
$ws.Range("I26").Value = 2.86

$ws.Range("C27").Value = 14590
$ws.Range("D27").Value = 0.0020999999999999999
$ws.Range("I27").Value = 4.46
$ws.Range("J27").Value = 87
$ws.Range("K27").Value = 87

$ws.Range("C28").Value = 13870
$ws.Range("D28").Value = 0
$ws.Range("D28").NumberFormat = "0%"
$ws.Range("I28").Value = 3.6
$ws.Range("J28").Value = 84
$ws.Range("K28").Value = 84

$ws.Range("C29").Value = 23100
$ws.Range("D29").Value = -0.017000000000000001
$ws.Range("I29").Value = 4.3099999999999996

$ws.Range("D30").Value = -0.012

